$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new year (2020) column has been appended as column O, mirroring the
# formatting already used for column N (the previous last year column).

# Header row (row 4): year label 2020
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("O4").Value = 2020

# Data row (row 5): value for 2020
$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("O5").Value = 83.3

$excel.CutCopyMode = $false

# Move/update the active selection like in the authored workbook
$ws.Range("O12").Select()
